# Apply the "Add files via upload" revision to the Testcases Report workbook.
#
# Net content changes (derived from the OOXML diff):
#   - Column E ("Pre-Requisite") gets populated for the four test-case rows
#     (6-9) with " User should have good internet connectivity."
#   - Row 6: the stray whitespace "Commnets" cell (K6) is cleared, and the
#     "TC for Automation(Y/N)" value (L6) flips from "Y" to "N".
#   - Row 7: Status (J7) flips from "Fail" to "Pass"; the "Commnets" (K7,
#     "Steps are not clear to follow") and "BUG ID" (M7, "BUG-1234") cells
#     are cleared out entirely.
#   - Rows 8 and 9 each gain a "N" in the "TC for Automation(Y/N)" column (L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Pre-Requisite" text for every test case row.
$prereq = " User should have good internet connectivity."
$ws.Range("E6").Value = $prereq
$ws.Range("E7").Value = $prereq
$ws.Range("E8").Value = $prereq
$ws.Range("E9").Value = $prereq

# Row 6 (HomePage_TC_OO1): clear the blank "Commnets" cell, flip automation flag.
$ws.Range("K6").Clear()
$ws.Range("L6").Value = "N"

# Row 7 (LoginPage_TC_OO2): Status Fail -> Pass, drop stale comment/bug id.
$ws.Range("J7").Value = "Pass"
$ws.Range("K7").Clear()
$ws.Range("M7").Clear()

# Rows 8 & 9 (LoginPage_TC_OO3 / Dashboard_TC_OO4): mark automation as "N".
$ws.Range("L8").Value = "N"
$ws.Range("L9").Value = "N"

# Match the author's final selection left in the saved workbook.
$ws.Range("M7").Select()
